$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-16 Thursday", "2023-11-17 Friday"),
    @("90×74=6660", "74×74=5476"),
    @("72×79=5688", "80×67=5360"),
    @("17×15=255", "83×18=1494"),
    @("51×34=1734", "22×54=1188"),
    @("75×27=2025", "15×35=525"),
    @("76×49=3724", "13×66=858"),
    @("52×59=3068", "22×74=1628"),
    @("88×96=8448", "96×40=3840"),
    @("90×87=7830", "44×62=2728"),
    @("70×68=4760", "81×35=2835"),
    @("74×76=5624", "80×93=7440"),
    @("18×93=1674", "67×22=1474"),
    @("35×27=945", "81×62=5022"),
    @("37×88=3256", "60×97=5820"),
    @("53×94=4982", "50×52=2600"),
    @("11×37=407", "37×94=3478"),
    @("47×93=4371", "69×53=3657"),
    @("16×32=512", "77×39=3003"),
    @("21×84=1764", "37×36=1332"),
    @("44×36=1584", "43×46=1978"),
    @("96×14=1344", "39×65=2535"),
    @("82×61=5002", "83×50=4150"),
    @("68×86=5848", "89×41=3649"),
    @("17×77=1309", "39×57=2223"),
    @("78×36=2808", "78×54=4212")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
